# "rerun benchmarks after removing irrelevant test"
# Updates the three benchmark rows in 'ubuntu-chrome' with the re-measured
# numbers, nudges the chart over to its new spot, and leaves the selection
# where the author's cursor ended up (H8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-run benchmark numbers -------------------------------------------------
$ws.Range("B2").Value = 11.3
$ws.Range("D2").Value = 34.8

$ws.Range("B3").Value = 11.2
$ws.Range("C3").Value = 10.3
$ws.Range("D3").Value = 35

$ws.Range("B4").Value = 11.2
$ws.Range("C4").Value = 10.2

# --- Reposition the chart slightly (dragged down/right a bit) ----------------
$co = $ws.ChartObjects().Item(1)
$co.Left = 356.49269685039405
$co.Top = 103.16535433070865
$co.Width = 414.70344488189
$co.Height = 255.21889763779495

# --- Restore the author's final selection -------------------------------------
$ws.Range("H8").Select()
